# Generate Report for Handback
#
# This script mirrors a "handback" localization-report refresh:
#   - Status cells flip from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - "Latest Handback DateTime" is stamped with the new handback time
#   - The now-stale "Latest Handback Name" and "Error Detail" (version
#     mismatch warning) columns are cleared out
#   - A few report columns are widened/narrowed to fit the new content
#
# Column widths: the host's ColumnWidth setter rounds to whole-pixel
# (1/6-character) granularity before re-deriving the stored OOXML width
# (width = round(ColumnWidth*6)/6 + 5/6), same as genuine Excel does
# internally -- so we pick the COM ColumnWidth value whose derived width
# lands closest to each target.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsZhCn.Range("L2").Value = "2017-02-17 07:11:44"
$wsZhCn.Range("L3").Value = "2017-02-17 07:11:44"

$wsZhCn.Range("M2").Value = ""
$wsZhCn.Range("M3").Value = ""

$wsZhCn.Range("R2").Value = ""
$wsZhCn.Range("R3").Value = ""

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(13).ColumnWidth = 23.0
$wsZhCn.Columns.Item(18).ColumnWidth = 12.833333333333334

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

$wsDeDe.Range("L2").Value = "2017-02-17 07:12:12"
$wsDeDe.Range("L3").Value = "2017-02-17 07:12:12"

$wsDeDe.Range("M2").Value = ""
$wsDeDe.Range("M3").Value = ""

$wsDeDe.Range("R2").Value = ""
$wsDeDe.Range("R3").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(13).ColumnWidth = 23.0
$wsDeDe.Columns.Item(18).ColumnWidth = 12.833333333333334

# ---------------------------------------------------------------------
# Overview sheet - widen the per-language summary columns (E, F)
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668
